$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-44 down to 42-45
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly entry
$ws.Range("A41").Value = 5
$ws.Range("B41").Value = "Macroferia Regional de Talca"
$ws.Range("C41").Value = "Maule"
$ws.Range("D41").Value = 44706
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100104
$ws.Range("H41").Value = "Frutos de pepita"
$ws.Range("I41").Value = 100104003
$ws.Range("J41").Value = "Membrillo"
$ws.Range("K41").Value = "Champion"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 300
$ws.Range("N41").Value = 10000
$ws.Range("O41").Value = 10000
$ws.Range("P41").Value = 10000
$ws.Range("Q41").Value = "$/caja 18 kilos granel"
$ws.Range("R41").Value = "Región de O'Higgins"
$ws.Range("S41").Value = 556
$ws.Range("T41").Value = 18
